$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# 1) Title paragraph: drop trailing space in "Software Requirement Specification "
$d.Content.Find.Execute("Software Requirement Specification ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Software Requirement Specification", 2)

# 2) New paragraph right after the title: "(SRS) " (two runs: "(SRS)" + " "), sz 36, centered
$pTitle = $d.Paragraphs(1)
$pTitle.Range.InsertParagraphAfter()
$pSrs = $d.Paragraphs(2)
$xmlSrs = '<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="center"/><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr>' + `
          '<w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>(SRS)</w:t></w:r>' + `
          '<w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>'
$pSrs.Range.InsertXML($xmlSrs)

# 3) Two new paragraphs after "There are four action actors..."
$pActors = $d.Paragraphs(3)

$pActors.Range.InsertParagraphAfter()
$pChat = $d.Paragraphs(4)
$chatText = "We need two type of chat, Private and Group. Private chat used to connect two member with some identity that it" + [char]0x2019 + "s isolated from others. However Group chat is like a room with an id where it can include more than two member and everybody joined, can see all Messages."
$xmlChat = '<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="left"/><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr>' + `
           '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>' + $chatText + '</w:t></w:r></w:p>'
$pChat.Range.InsertXML($xmlChat)

$pChat2 = $d.Paragraphs(4)
$pChat2.Range.InsertParagraphAfter()
$pBlank = $d.Paragraphs(5)
$xmlBlank = '<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="left"/><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr>' + `
            '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>'
$pBlank.Range.InsertXML($xmlBlank)
# InsertXML drops <w:spacing>; set it through the paragraph-format properties instead.
$pBlank2 = $d.Paragraphs(5)
$pBlank2.Format.SpaceBefore = 0
$pBlank2.Format.SpaceAfter = 8

Write-Output $d.Paragraphs.Count
foreach ($p in $d.Paragraphs) {
    Write-Output ("[" + $p.Range.Text + "]")
}
